# Scheduled market-data refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit figures (columns H-N) on the affected leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Values come straight from the upstream market-board snapshot; cells that
# the snapshot no longer populates are cleared (set to $null) so they stay
# absent from the sheet, matching cells that newly appear are created by
# simply assigning a value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null

$ws.Range("H28").Value = 805.1667
$ws.Range("I28").Value = 757.75
$ws.Range("J28").Value = 900
$ws.Range("K28").Value = 757.75
$ws.Range("L28").Value = 900
$ws.Range("M28").Value = -272.75
$ws.Range("N28").Value = -1870

$ws.Range("H58").Value = 1698.5186
$ws.Range("I58").Value = 564.61536
$ws.Range("J58").Value = 2751.4285
$ws.Range("K58").Value = 1693.84608
$ws.Range("L58").Value = 8254.2855
$ws.Range("M58").Value = -1543.84608
$ws.Range("N58").Value = -8554.2855

$ws.Range("H80").Value = 583.375
$ws.Range("I80").Value = 493.4
$ws.Range("J80").Value = 733.3333
$ws.Range("K80").Value = 1480.2
$ws.Range("L80").Value = 2199.9999
$ws.Range("M80").Value = -482.1999999999998
$ws.Range("N80").Value = -4195.9999

$ws.Range("H83").Value = 583.375
$ws.Range("I83").Value = 493.4
$ws.Range("J83").Value = 733.3333
$ws.Range("K83").Value = 4440.599999999999
$ws.Range("L83").Value = 6599.9997
$ws.Range("M83").Value = 551.4000000000005
$ws.Range("N83").Value = -16583.9997

$ws.Range("H88").Value = 1126.5883
$ws.Range("I88").Value = 889.8
$ws.Range("J88").Value = 1225.25
$ws.Range("K88").Value = 889.8
$ws.Range("L88").Value = 1225.25
$ws.Range("M88").Value = -483.8
$ws.Range("N88").Value = -2037.25

$ws.Range("H91").Value = 1126.5883
$ws.Range("I91").Value = 889.8
$ws.Range("J91").Value = 1225.25
$ws.Range("K91").Value = 889.8
$ws.Range("L91").Value = 1225.25
$ws.Range("M91").Value = 514.2
$ws.Range("N91").Value = -4033.25

$ws.Range("H92").Value = 280.3684
$ws.Range("I92").Value = 175.6
$ws.Range("K92").Value = 175.6
$ws.Range("M92").Value = 1072.4

$ws.Range("H94").Value = 7846.9165
$ws.Range("I94").Value = 7846.9165
$ws.Range("K94").Value = 7846.9165
$ws.Range("M94").Value = -7395.9165

$ws.Range("H97").Value = 2889.1428
$ws.Range("J97").Value = 2889.1428
$ws.Range("L97").Value = 8667.428400000001
$ws.Range("N97").Value = -9659.428400000001

$ws.Range("H98").Value = 1340.6875
$ws.Range("I98").Value = 1303.7858
$ws.Range("K98").Value = 1303.7858
$ws.Range("M98").Value = 194.2141999999999

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null

$ws.Range("H107").Value = 1225.3914
$ws.Range("I107").Value = 1371.3889
$ws.Range("K107").Value = 1371.3889
$ws.Range("M107").Value = 548.6111000000001

$ws.Range("H112").Value = 5703.3335
$ws.Range("I112").Value = 3500
$ws.Range("J112").Value = 10110
$ws.Range("K112").Value = 10500
$ws.Range("L112").Value = 30330
$ws.Range("M112").Value = -9392
$ws.Range("N112").Value = -32546

$ws.Range("H122").Value = 1340.6875
$ws.Range("I122").Value = 1303.7858
$ws.Range("K122").Value = 3911.3574
$ws.Range("M122").Value = -1461.3574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1691.3334
$ws.Range("J74").Value = 2031
$ws.Range("L74").Value = 2031
$ws.Range("N74").Value = -3779

$ws.Range("H77").Value = 1691.3334
$ws.Range("J77").Value = 2031
$ws.Range("L77").Value = 10155
$ws.Range("N77").Value = -18891

$ws.Range("H102").Value = 10106158
$ws.Range("I102").Value = 11116574
$ws.Range("K102").Value = 11116574
$ws.Range("M102").Value = -11114952

$ws.Range("H122").Value = 12176.333
$ws.Range("I122").Value = 8790.440000000001
$ws.Range("J122").Value = 54500
$ws.Range("K122").Value = 26371.32
$ws.Range("L122").Value = 163500
$ws.Range("M122").Value = -23921.32
$ws.Range("N122").Value = -168400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2673.4546
$ws.Range("I105").Value = 2015.8518
$ws.Range("K105").Value = 2015.8518
$ws.Range("M105").Value = -268.8517999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 170
$ws.Range("J15").Value = 275
$ws.Range("L15").Value = 275
$ws.Range("N15").Value = -615

$ws.Range("H22").Value = 89273.336
$ws.Range("I22").Value = 132844.33
$ws.Range("J22").Value = 23916.834
$ws.Range("K22").Value = 132844.33
$ws.Range("L22").Value = 23916.834
$ws.Range("M22").Value = -132494.33
$ws.Range("N22").Value = -24616.834

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1199
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = $null

$ws.Range("H63").Value = 3100
$ws.Range("I63").Value = 3100
$ws.Range("K63").Value = 9300
$ws.Range("M63").Value = -8551

$ws.Range("H66").Value = 3100
$ws.Range("I66").Value = 3100
$ws.Range("K66").Value = 27900
$ws.Range("M66").Value = -24156

$ws.Range("H92").Value = 900
$ws.Range("J92").Value = 900
$ws.Range("L92").Value = 2700
$ws.Range("N92").Value = -5196

$ws.Range("H106").Value = 8666.666999999999
$ws.Range("J106").Value = 8666.666999999999
$ws.Range("L106").Value = 26000.001
$ws.Range("N106").Value = -27892.001

$ws.Range("H109").Value = 1127.8572
$ws.Range("I109").Value = 816
$ws.Range("K109").Value = 2448
$ws.Range("M109").Value = -1408

$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 2000
$ws.Range("K125").Value = 6000
$ws.Range("M125").Value = -1080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 343.6111
$ws.Range("I2").Value = 95
$ws.Range("J2").Value = 501.81818
$ws.Range("K2").Value = 95
$ws.Range("L2").Value = 501.81818
$ws.Range("M2").Value = 18
$ws.Range("N2").Value = -727.81818

$ws.Range("H113").Value = 13916.667
$ws.Range("I113").Value = 1250
$ws.Range("K113").Value = 1250
$ws.Range("M113").Value = 920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3199.75
$ws.Range("I100").Value = 3166.3333
$ws.Range("K100").Value = 3166.3333
$ws.Range("M100").Value = -2625.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1705.4546
$ws.Range("I100").Value = 1935.6666
$ws.Range("K100").Value = 3871.3332
$ws.Range("M100").Value = -3330.3332

$ws.Range("H122").Value = 3379.8235
$ws.Range("I122").Value = 3042.5454
$ws.Range("J122").Value = 3998.1667
$ws.Range("K122").Value = 9127.636200000001
$ws.Range("L122").Value = 11994.5001
$ws.Range("M122").Value = -6677.636200000001
$ws.Range("N122").Value = -16894.5001
